$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Rewrite the "Plot ... the geostrophic wind speed?" sentence into the
#    new, longer sentence, split across multiple runs exactly as in the
#    target revision (the runs all share the same rPr, so we deliberately
#    flip a formatting property on/off on each chunk to force Word to keep
#    them as separate <w:r> elements instead of re-coalescing them).
# ---------------------------------------------------------------------
$p3 = $d.Paragraphs(3)
$oldSentence = " the angular deviation of the surface wind relative to the geostrophic wind as a function of the geostrophic wind speed? "

$findRng = $p3.Range.Duplicate
$found = $findRng.Find.Execute($oldSentence, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $startPos = $findRng.Start

    $newText = " the angular deviation of the surface wind relative to the geostrophic wind  (the cross-isobar flow angle) as a function of RO = G/(f * z_0), as in Garratt (1992) figure 12b? Or just use the geostrophic wind speed, as f * z_0 is constant. "
    $findRng.Text = $newText

    # Lengths (in characters) of each run-chunk of $newText, in order:
    #  " the angular deviation of the surface wind relative to the geostrophic wind "
    #  " ("
    #  "the cross-isobar flow angle) "
    #  "as a function of"
    #  " RO = G/(f "
    #  "* "
    #  "z_0), as in Garratt (1992) figure 12b?"
    #  " Or just use the "
    #  "geostrophic wind speed, as f * z_0 is constant. "
    $chunkLens = @(76, 2, 29, 16, 11, 2, 38, 17, 48)

    $pos = $startPos
    foreach ($len in $chunkLens) {
        $chunkRng = $d.Range($pos, $pos + $len)
        $chunkRng.Font.Bold = $true
        $chunkRng.Font.Bold = $false
        $pos = $pos + $len
    }
}

# ---------------------------------------------------------------------
# 2. Move the "_GoBack" bookmark: delete it from its old location (after
#    "Look at occurrence of inertial oscillations") and re-add it inside
#    the sentence just edited, between "* " and "z_0), as in Garratt...".
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$p3b = $d.Paragraphs(3)
$bmRng = $p3b.Range.Duplicate
$bmRng.Find.Execute("z_0), as in Garratt (1992) figure 12b?", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bmPos = $bmRng.Start
$d.Bookmarks.Add("_GoBack", $d.Range($bmPos, $bmPos))

# ---------------------------------------------------------------------
# 3. Add the "Nadruk" (Emphasis) character style to styles.xml.
# ---------------------------------------------------------------------
$st = $d.Styles.Add("Nadruk", 2)
$st.QuickStyle = $true
$st.Priority = 20
$st.BaseStyle = "Standaardalinea-lettertype"
$st.NameLocal = "Emphasis"
$st.Font.Italic = $true
$st.Font.ItalicBi = $true
